$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Time Left" column (K) ---
# Header K1: set text, then copy direct formatting (font/border/alignment)
# from the neighboring header cell J1 so it matches the rest of the header row.
$ws.Range('K1').Value = 'Time Left'
$ws.Range('J1').Copy()
$ws.Range('K1').PasteSpecial(-4122)

# --- Refresh eBay listing URLs in column C (tracking/session params updated) ---
$ws.Range('C2').Value = 'https://www.ebay.com/itm/225299351425?hash=item3474e34781:g:WagAAOSwxeNjmfpZ&amdata=enc%3AAQAHAAAAoMZmeAwO9uGAOepXocLGZPSvZ7HVASv%2Bz9HfSuZcidVwulsj3Q3WtshmmFTonlSzq5%2BY9bf%2BFpgGcY40jugzJBklU60ZTfSRPH1tr4%2FrOVDJ9glCe1RFfTe5D7QPj4U00ElLi3rIDkJRzuDDZaFAQNrqXetgCj2i6V5gt%2FU9Nq%2BDkEp%2F12EGRZCz%2FiXftz2q%2BuTdoaFyU%2BG3PUtDgiKh0vo%3D%7Ctkp%3ABk9SR9bSt7WjYQ'
$ws.Range('C3').Value = 'https://www.ebay.com/itm/404050551629?epid=8047515616&hash=item5e134a234d:g:cYsAAOSwWPljk3xw&amdata=enc%3AAQAHAAAAoGiTEUPW5VtRzKGDLjO84d5VQscWIRHpASDiRaUr3QjG4U1Tfd8mHVtqMu7m01S817jmRLNDFxk3hWr3zQjbO73QXJmWfucg84ep%2F4517NjXyAbmPV2sr97eTBUQVnyaKosudO8s7b5bnhjEEp2QALiD7CpYiAR6kf0IafFKuFvKeGJUws3gzUVGOJblJWiiNsPxSqlykHIzR0%2F0RWlr1WA%3D%7Ctkp%3ABk9SR9bSt7WjYQ'
$ws.Range('C4').Value = 'https://www.ebay.com/itm/295403696780?epid=7050209760&hash=item44c76eb28c:g:bw4AAOSwChJjk1Kw&amdata=enc%3AAQAHAAAAkBAKU7HgZCqFBzuU0e3plvMFS2not8HtJB9dLimFOyJnFvxXxCxSMCGsQ141YG1jpLb98QSxAQrNYtKTjol6U8fjEjTcp6xDUIPHuSB%2F0X274OLeRxUcpTF24YDK0Wheqj%2B7UMxGoYlAQNTGE3CVsB2f5J%2FQEq%2Bz6X74BusFiv5gXqJo3hZR0oO%2BZ1vcI11c1Q%3D%3D%7Ctkp%3ABk9SR9bSt7WjYQ'
$ws.Range('C5').Value = 'https://www.ebay.com/itm/295418574962?hash=item44c851b872:g:QQkAAOSwJNhjmuak&amdata=enc%3AAQAHAAAAoInV0dSL5sn815Yx6uFdOd3lkWn1srcTHg1wPp96pvDm9zAsg8a1xELtPQViMrGfgC5Mz5WhCmnSSbqgUrAptcBoIn5VHUaK90CSTzEEhr369OFFur9mEm7yI8etIDRA9C%2FVHdfAVvIL8i41Uumy7s7htLe7oYC1sXy4BwQ4hARbN%2FMQJ8Tm5zyJNhzt2rkyeqgFfq78NNFeGc3v8qDueOo%3D%7Ctkp%3ABk9SR9bSt7WjYQ'
$ws.Range('C6').Value = 'https://www.ebay.com/itm/285074976656?epid=14050520857&hash=item425fcaef90:g:EGcAAOSwdJBjlNZH&amdata=enc%3AAQAHAAAAoJj3syuSee2Bl6pg78DOlpmggJb6JULGmdUQt%2FVrBsLHD3FFtxjIlg0ihGLMNLeO%2F6pEDHAPAgl1GlQ46bvqTy9azvTM%2B2nvBJ6Lqr7YQqckMY%2FbkEgloFUJbwM5gvriIzHW9rq7x6luVhkX0UIOtL822RfKlxIVMXaw9wQjvrF9qN0dbmgVLpjf9NfQLalpPmiScXi15gy4CUktr9NfRsA%3D%7Ctkp%3ABk9SR9bSt7WjYQ'

# --- Populate "Time Left" values for each auction row ---
$ws.Range('K2').Value = 'Not specified'
$ws.Range('K3').Value = 'Not specified'
$ws.Range('K4').Value = '3d restantes'
$ws.Range('K5').Value = '5d 19h restantes'
$ws.Range('K6').Value = '2d 4h restantes'

# --- Remove stale "Novo anúncio" (New listing) prefix from row 5 title ---
$ws.Range('B5').Value = 'Tarjeta Gráfica AMD Radeon Gigabyte RX 6900 XT Gaming OC 16 GB'

